$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# Locate the paragraph whose text is exactly "lol" (the one being edited per the diff).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "lol") {
        $target = $p
    }
}

if ($target -ne $null) {
    $r = $target.Range
    $paraStart = $r.Start
    $paraEnd = $r.End

    # Replace the run text "lol" -> two runs "L" + "ol", leaving the trailing
    # paragraph mark (at paraEnd-1..paraEnd) untouched.
    $textRange = $d.Range($paraStart, $paraEnd - 1)
    $textRange.InsertXML("<w:p xmlns:w='$wNs'><w:r><w:t>L</w:t></w:r><w:r><w:t>ol</w:t></w:r></w:p>")

    # Re-acquire the (now "Lol") paragraph end, then append: blank paragraph,
    # "V3" paragraph, blank paragraph, right after it.
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.TrimEnd("`r`a") -eq "Lol") {
            $insertAt = $p.Range.End
            $afterRange = $d.Range($insertAt, $insertAt)
            $afterRange.InsertXML("<w:p xmlns:w='$wNs'/><w:p xmlns:w='$wNs'><w:r><w:t>V3</w:t></w:r></w:p><w:p xmlns:w='$wNs'/>")
        }
    }
}
